# Clarified the input specification
# Update the "Comment/Description" column (column E) text for several
# variables in the "Raw input data" sheet to add clarifying notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# KEY_DONOR_DOB description: clarify the date format
$ws.Range("E5").Value = "Date of birth (e.g. 19580111, yyyymmdd)"

# FERRITIN_FIRST description: mark as optional
$ws.Range("E7").Value = "OPTIONAL, First ferritin (ug/L = ng/mL)"

# FERRITIN_LAST description: mark as optional
$ws.Range("E8").Value = "OPTIONAL, Last ferritin (ug/L = ng/mL)"

# FERRITIN_LAST_DATE description: mark as optional
$ws.Range("E9").Value = "OPTIONAL, Date when the last ferritin was measured (yyyymmdd)"

# DONAT_STATUS description: expand on possible status codes and usage note
$ws.Range("E13").Value = "Donation status (- =OK, V=change of the bag, R=disposable blood unit, E=no donation, D=does not fulfill requirements, …). I only check whether status is equal to ‘-’ or not."

# KEY_DONAT_PHLEB description: expand on possible donation type codes and usage note
$ws.Range("E14").Value = "Donation Type (K=Whole Blood donation, P=Plasmapheresis, T=Trombapheresis, =No donation, H=Whole blood – not studied,…). I only compare whether PHLEB == or != ‘K’, ‘H’, ‘*’"

# Restore the active selection to match the author's final cursor position
[void]$ws.Range("E18").Select()
